{"js": "// The author removed two to-do bullets (\"\u0421\u0434\u0435\u043b\u0430\u0442\u044c \u0431\u043b\u043e\u043a\u0438 \u043f\u0443\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044f...\" and\n// \"\u0423\u0441\u0442\u0430\u0432\u043a\u0438 \u0434\u043b\u044f \u043f\u0438\u0434\u0430...\") and removed the last bullet (\"\u041f\u0440\u0435\u0434\u0443\u0441\u043c\u043e\u0442\u0440\u0435\u0442\u044c \u043f\u0440\u0438\n// \u043e\u043a\u043e\u043d\u0447\u0430\u043d\u0438\u0438 \u0441\u0442\u0435\u0440\u0438\u043b\u0438\u0437\u0430\u0446\u0438\u0438...\") that carried the `_GoBack` bookmark. That\n// bookmark now opens the \"\u0411\u043b\u043e\u043a\u0438 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043d\u0430\u0441\u043e\u0441\u0430\u043c\u0438 \u2013 \u041d\u0430\u0437\u043d\u0430\u0447\u0438\u0442\u044c:\" bullet\n// instead, so the editing caret returns there on reopen.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet pSdelatBloki = null;      // \"\u0421\u0434\u0435\u043b\u0430\u0442\u044c \u0431\u043b\u043e\u043a\u0438 \u043f\u0443\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044f \u043a\u043b\u0430\u043f\u0430\u043d\u0430\u043c\u0438...\"\nlet pUstavki = null;          // \"\u0423\u0441\u0442\u0430\u0432\u043a\u0438 \u0434\u043b\u044f \u043f\u0438\u0434\u0430 \u043f\u0440\u043e\u043f\u0438\u0441\u0430\u0442\u044c...\"\nlet pNasosyNaznachit = null;  // \"\u0411\u043b\u043e\u043a\u0438 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043d\u0430\u0441\u043e\u0441\u0430\u043c\u0438 \u2013 \u041d\u0430\u0437\u043d\u0430\u0447\u0438\u0442\u044c:\"\nlet pPredusmotret = null;     // \"\u041f\u0440\u0435\u0434\u0443\u0441\u043c\u043e\u0442\u0440\u0435\u0442\u044c \u043f\u0440\u0438 \u043e\u043a\u043e\u043d\u0447\u0430\u043d\u0438\u0438 \u0441\u0442\u0435\u0440\u0438\u043b\u0438\u0437\u0430\u0446\u0438\u0438...\"\n\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  if (text.indexOf(\"\u0421\u0434\u0435\u043b\u0430\u0442\u044c \u0431\u043b\u043e\u043a\u0438\") === 0) {\n    pSdelatBloki = p;\n  } else if (text.indexOf(\"\u0423\u0441\u0442\u0430\u0432\u043a\u0438\") === 0) {\n    pUstavki = p;\n  } else if (text.indexOf(\"\u0411\u043b\u043e\u043a\u0438 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043d\u0430\u0441\u043e\u0441\u0430\u043c\u0438\") === 0) {\n    pNasosyNaznachit = p;\n  } else if (text.indexOf(\"\u041f\u0440\u0435\u0434\u0443\u0441\u043c\u043e\u0442\u0440\u0435\u0442\u044c \u043f\u0440\u0438 \u043e\u043a\u043e\u043d\u0447\u0430\u043d\u0438\u0438\") === 0) {\n    pPredusmotret = p;\n  }\n}\n\nif (!pSdelatBloki || !pUstavki || !pNasosyNaznachit || !pPredusmotret) {\n  throw new Error(\"Could not locate one of the expected paragraphs\");\n}\n\n// Drop the two obsolete task bullets entirely.\npSdelatBloki.delete();\npUstavki.delete();\n\n// Drop the final bullet too \u2014 this also removes the `_GoBack` bookmark it\n// used to carry (the bookmark is fully contained inside that paragraph).\npPredusmotret.delete();\n\n// Re-create the `_GoBack` bookmark at the very start of the\n// \"\u0411\u043b\u043e\u043a\u0438 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043d\u0430\u0441\u043e\u0441\u0430\u043c\u0438 \u2013 \u041d\u0430\u0437\u043d\u0430\u0447\u0438\u0442\u044c:\" paragraph.\npNasosyNaznachit.getRange(\"Start\").insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The author removed two to-do bullets (\"\u0421\u0434\u0435\u043b\u0430\u0442\u044c \u0431\u043b\u043e\u043a\u0438 \u043f\u0443\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044f...\" and\n# \"\u0423\u0441\u0442\u0430\u0432\u043a\u0438 \u0434\u043b\u044f \u043f\u0438\u0434\u0430...\") and removed the last bullet (\"\u041f\u0440\u0435\u0434\u0443\u0441\u043c\u043e\u0442\u0440\u0435\u0442\u044c \u043f\u0440\u0438\n# \u043e\u043a\u043e\u043d\u0447\u0430\u043d\u0438\u0438 \u0441\u0442\u0435\u0440\u0438\u043b\u0438\u0437\u0430\u0446\u0438\u0438...\") that carried the `_GoBack` bookmark. That\n# bookmark now opens the \"\u0411\u043b\u043e\u043a\u0438 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043d\u0430\u0441\u043e\u0441\u0430\u043c\u0438 \u2013 \u041d\u0430\u0437\u043d\u0430\u0447\u0438\u0442\u044c:\" bullet\n# instead, so the editing caret returns there on reopen.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphStartingWith($prefix) {\n  foreach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"$prefix*\") {\n      return $p\n    }\n  }\n  return $null\n}\n\n# Delete from the bottom of the document upward so previously-found\n# paragraph objects that sit earlier in the document stay valid. (Re-locate\n# each paragraph right before acting on it, rather than caching stale\n# references, so a prior delete can't leave us pointing at the wrong text.)\n\n$pPredusmotret = Find-ParagraphStartingWith(\"\u041f\u0440\u0435\u0434\u0443\u0441\u043c\u043e\u0442\u0440\u0435\u0442\u044c \u043f\u0440\u0438 \u043e\u043a\u043e\u043d\u0447\u0430\u043d\u0438\u0438\")\n$pPredusmotret.Range.Delete()\n\n$pUstavki = Find-ParagraphStartingWith(\"\u0423\u0441\u0442\u0430\u0432\u043a\u0438\")\n$pUstavki.Range.Delete()\n\n$pSdelatBloki = Find-ParagraphStartingWith(\"\u0421\u0434\u0435\u043b\u0430\u0442\u044c \u0431\u043b\u043e\u043a\u0438\")\n$pSdelatBloki.Range.Delete()\n\n# Re-create the `_GoBack` bookmark at the very start of the\n# \"\u0411\u043b\u043e\u043a\u0438 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043d\u0430\u0441\u043e\u0441\u0430\u043c\u0438 \u2013 \u041d\u0430\u0437\u043d\u0430\u0447\u0438\u0442\u044c:\" paragraph.\n$pNasosyNaznachit = Find-ParagraphStartingWith(\"\u0411\u043b\u043e\u043a\u0438 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043d\u0430\u0441\u043e\u0441\u0430\u043c\u0438\")\n$r = $pNasosyNaznachit.Range.Duplicate()\n$r.Collapse(1)   # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $r)\n"}
